# The presentation's theme (ppt/theme/theme1.xml, used by the slide
# master / all slides) is switched from the "Integral" palette to the
# stock "Office Theme" palette. dk1/lt1 (black/white) are unchanged in
# both palettes; dk2, lt2 and the six accents plus the two hyperlink
# colors change.
#
# Font scheme (majorFont/minorFont) and format scheme (fill/line/effect
# styles) are already identical between the old and new theme, so only
# the color scheme needs to change.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72
